# Update "想去人数" (want-to-go count) and a couple of "最低票价" (min price)
# values across the 展览 (Exhibition), 演出 (Show), and 全部类型 (All types)
# sheets, matching refreshed scrape output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 91
$ws.Range("F4").Value = 9909
$ws.Range("F5").Value = 692
$ws.Range("F6").Value = 181
$ws.Range("F7").Value = 364
$ws.Range("F8").Value = 396
$ws.Range("F9").Value = 443
$ws.Range("F11").Value = 220
$ws.Range("F13").Value = 498
$ws.Range("F14").Value = 12538
$ws.Range("F15").Value = 41
$ws.Range("F19").Value = 258
$ws.Range("F22").Value = 135
$ws.Range("F23").Value = 171
$ws.Range("F26").Value = 91
$ws.Range("F28").Value = 65
$ws.Range("F29").Value = 2164
$ws.Range("F30").Value = 1068
$ws.Range("F31").Value = 4243
$ws.Range("F32").Value = 3761
$ws.Range("F33").Value = 775
$ws.Range("F35").Value = 3076
$ws.Range("F36").Value = 55
$ws.Range("F37").Value = 1357
$ws.Range("F38").Value = 207
$ws.Range("F40").Value = 39
$ws.Range("F41").Value = 127
$ws.Range("F42").Value = 472
$ws.Range("F43").Value = 620
$ws.Range("F44").Value = 74
$ws.Range("F46").Value = 275
$ws.Range("F49").Value = 160

# ---- Sheet: 演出 (Show) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G5").Value = 280
$ws.Range("F11").Value = 31
$ws.Range("F22").Value = 38

# ---- Sheet: 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 91
$ws.Range("F6").Value = 9909
$ws.Range("F7").Value = 692
$ws.Range("G8").Value = 280
$ws.Range("F9").Value = 181
$ws.Range("F10").Value = 364
$ws.Range("F11").Value = 396
$ws.Range("F12").Value = 443
$ws.Range("F14").Value = 220
$ws.Range("F15").Value = 498
$ws.Range("F16").Value = 12538
$ws.Range("F19").Value = 258
$ws.Range("F22").Value = 135
$ws.Range("F23").Value = 171
$ws.Range("F26").Value = 91
$ws.Range("F27").Value = 65
$ws.Range("F28").Value = 2164
$ws.Range("F29").Value = 1068
$ws.Range("F30").Value = 4243
$ws.Range("F31").Value = 3761
$ws.Range("F32").Value = 775
$ws.Range("F34").Value = 3076
$ws.Range("F35").Value = 55
$ws.Range("F36").Value = 1357
$ws.Range("F37").Value = 207
$ws.Range("F39").Value = 39
$ws.Range("F40").Value = 127
$ws.Range("F41").Value = 472
$ws.Range("F42").Value = 38
$ws.Range("F43").Value = 620
$ws.Range("F44").Value = 74
$ws.Range("F46").Value = 275
$ws.Range("F49").Value = 160
